$wb = $excel.ActiveWorkbook

# "建物" (Building) sheet: property_category column (I) was incorrectly set to
# "land" for all rows; correct it to "building".
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2:I8").Value = "building"

# "汽車" (Car) sheet: property_category column (H) was incorrectly set to
# "land" for all rows; correct it to "car".
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2:H3").Value = "car"
